{"js": "// Remove the opening paragraph of the document (the \"This book will help you\n// master R plots...\" paragraph, which includes the Leanpub hyperlink and a\n// line break), leaving the document starting at \"We have spent a long time...\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items.find(p =>\n  p.text.indexOf(\"This book will help you master R plots the easy way\") !== -1\n);\n\nif (target) {\n  target.delete();\n} else {\n  // Fallback: if the expected text isn't found (e.g. already removed),\n  // do nothing further.\n  paragraphs.items[0].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the opening paragraph of the document (the \"This book will help you\n# master R plots...\" paragraph, which includes the Leanpub hyperlink and a\n# line break), leaving the document starting at \"We have spent a long time...\".\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*This book will help you master R plots the easy way*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $target.Range.Delete()\n} else {\n    # Fallback: remove the first paragraph if the expected text isn't found.\n    $d.Paragraphs.Item(1).Range.Delete()\n}\n"}
